$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1082.3846
$ws.Range("I33").Value = 1164.25
$ws.Range("K33").Value = 1164.25
$ws.Range("M33").Value = -935.25
$ws.Range("H70").Value = 6195.2246
$ws.Range("I70").Value = 4599.5
$ws.Range("K70").Value = 13798.5
$ws.Range("M70").Value = -13528.5
$ws.Range("H73").Value = 6195.2246
$ws.Range("I73").Value = 4599.5
$ws.Range("K73").Value = 13798.5
$ws.Range("M73").Value = -12862.5
$ws.Range("H86").Value = 2798
$ws.Range("J86").Value = 3259.1538
$ws.Range("L86").Value = 3259.1538
$ws.Range("N86").Value = -5505.1538
$ws.Range("H89").Value = 2798
$ws.Range("J89").Value = 3259.1538
$ws.Range("L89").Value = 16295.769
$ws.Range("N89").Value = -27527.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4956.2607
$ws.Range("I32").Value = 3504.0205
$ws.Range("K32").Value = 3504.0205
$ws.Range("M32").Value = -3217.0205
$ws.Range("H97").Value = 13091.913
$ws.Range("I97").Value = 10163.474
$ws.Range("K97").Value = 10163.474
$ws.Range("M97").Value = -9667.474
$ws.Range("H110").Value = 766.19354
$ws.Range("I110").Value = 831.9167
$ws.Range("J110").Value = 540.8570999999999
$ws.Range("K110").Value = 831.9167
$ws.Range("L110").Value = 540.8570999999999
$ws.Range("M110").Value = 1213.0833
$ws.Range("N110").Value = -4630.8571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4437.3076
$ws.Range("I86").Value = 5524.0435
$ws.Range("J86").Value = 2875.125
$ws.Range("K86").Value = 5524.0435
$ws.Range("L86").Value = 2875.125
$ws.Range("M86").Value = -4401.0435
$ws.Range("N86").Value = -5121.125
$ws.Range("H89").Value = 4437.3076
$ws.Range("I89").Value = 5524.0435
$ws.Range("J89").Value = 2875.125
$ws.Range("K89").Value = 27620.2175
$ws.Range("L89").Value = 14375.625
$ws.Range("M89").Value = -22004.2175
$ws.Range("N89").Value = -25607.625
$ws.Range("H94").Value = 3135.0667
$ws.Range("I94").Value = 506.88
$ws.Range("K94").Value = 506.88
$ws.Range("M94").Value = -55.88
$ws.Range("H107").Value = 2195
$ws.Range("J107").Value = 138
$ws.Range("L107").Value = 138
$ws.Range("N107").Value = -3978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2615.3125
$ws.Range("I16").Value = 2077.5715
$ws.Range("J16").Value = 3033.5557
$ws.Range("K16").Value = 2077.5715
$ws.Range("L16").Value = 3033.5557
$ws.Range("M16").Value = -1790.5715
$ws.Range("N16").Value = -3607.5557
$ws.Range("H58").Value = 4376.3477
$ws.Range("I58").Value = 4356.421
$ws.Range("J58").Value = 4471
$ws.Range("K58").Value = 4356.421
$ws.Range("L58").Value = 4471
$ws.Range("M58").Value = -4153.421
$ws.Range("N58").Value = -4877
$ws.Range("H99").Value = 4408.636
$ws.Range("I99").Value = 3583.3333
$ws.Range("J99").Value = 5399
$ws.Range("K99").Value = 3583.3333
$ws.Range("L99").Value = 5399
$ws.Range("M99").Value = -2085.3333
$ws.Range("N99").Value = -8395
$ws.Range("H113").Value = 2615.3125
$ws.Range("I113").Value = 2077.5715
$ws.Range("J113").Value = 3033.5557
$ws.Range("K113").Value = 2077.5715
$ws.Range("L113").Value = 3033.5557
$ws.Range("M113").Value = 92.42849999999999
$ws.Range("N113").Value = -7373.5557
$ws.Range("H114").Value = 37500
$ws.Range("J114").Value = 37500
$ws.Range("L114").Value = 37500
$ws.Range("N114").Value = -46178
$ws.Range("H117").Value = 84500
$ws.Range("J117").Value = 84500
$ws.Range("L117").Value = 84500
$ws.Range("N117").Value = -93678
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H121").Value = 100000
$ws.Range("J121").Value = 100000
$ws.Range("L121").Value = 100000
$ws.Range("N121").Value = -102620
$ws.Range("H126").Value = 4408.636
$ws.Range("I126").Value = 3583.3333
$ws.Range("J126").Value = 5399
$ws.Range("K126").Value = 10749.9999
$ws.Range("L126").Value = 16197
$ws.Range("M126").Value = -8279.999899999999
$ws.Range("N126").Value = -21137
$ws.Range("H136").Value = 4376.3477
$ws.Range("I136").Value = 4356.421
$ws.Range("J136").Value = 4471
$ws.Range("K136").Value = 13069.263
$ws.Range("L136").Value = 13413
$ws.Range("M136").Value = -10519.263
$ws.Range("N136").Value = -18513

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3500.625
$ws.Range("I138").Value = 3500.625
$ws.Range("K138").Value = 10501.875
$ws.Range("M138").Value = -5361.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 45000
$ws.Range("J110").Value = 45000
$ws.Range("L110").Value = 45000
$ws.Range("N110").Value = -53180
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H119").Value = 91000
$ws.Range("J119").Value = 91000
$ws.Range("L119").Value = 91000
$ws.Range("N119").Value = -100676
$ws.Range("H132").Value = 2839.5715
$ws.Range("I132").Value = 2283.147
$ws.Range("J132").Value = 4100.8
$ws.Range("K132").Value = 6849.441
$ws.Range("L132").Value = 12302.4
$ws.Range("M132").Value = -4319.441
$ws.Range("N132").Value = -17362.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 35500
$ws.Range("I74").Value = 35500
$ws.Range("K74").Value = 35500
$ws.Range("M74").Value = -34502
$ws.Range("H77").Value = 35500
$ws.Range("I77").Value = 35500
$ws.Range("K77").Value = 106500
$ws.Range("M77").Value = -101508
$ws.Range("H120").Value = 10000
$ws.Range("J120").Value = 10000
$ws.Range("L120").Value = 10000
$ws.Range("N120").Value = -19676

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6144.8354
$ws.Range("I62").Value = 3157.5757
$ws.Range("J62").Value = 8609.325000000001
$ws.Range("K62").Value = 3157.5757
$ws.Range("L62").Value = 8609.325000000001
$ws.Range("M62").Value = -2533.5757
$ws.Range("N62").Value = -9857.325000000001
$ws.Range("H65").Value = 6144.8354
$ws.Range("I65").Value = 3157.5757
$ws.Range("J65").Value = 8609.325000000001
$ws.Range("K65").Value = 15787.8785
$ws.Range("L65").Value = 43046.625
$ws.Range("M65").Value = -12667.8785
$ws.Range("N65").Value = -49286.625
$ws.Range("H126").Value = 1657.0333
$ws.Range("I126").Value = 1556.64
$ws.Range("K126").Value = 4669.92
$ws.Range("M126").Value = -2199.92

Write-Host "Applied all changes"
